$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '21.776.11'
$ws.Range('E2').Value = '  -1.71%  '
$ws.Range('D3').Value = '1.541.79'
$ws.Range('E3').Value = '  -1.33%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('E5').Value = '  +0.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '289.61'
$ws.Range('E6').Value = '  -0.13%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3894'
$ws.Range('E7').Value = '  +2.64%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3193'
$ws.Range('E8').Value = '  -2.90%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '43.45'
$ws.Range('E9').Value = '  -0.52%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07197'
$ws.Range('E10').Value = '  -2.39%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.065'
$ws.Range('E11').Value = '  -6.67%  '
$ws.Range('E12').Value = '  +0.02%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.628'
$ws.Range('E13').Value = '  -3.55%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.59'
$ws.Range('E14').Value = '  -7.14%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.599'
$ws.Range('E15').Value = '  -4.15%  '
$ws.Range('D16').Value = '1.542.58'
$ws.Range('E16').Value = '  -1.54%  '
$ws.Range('E17').Value = '  +1.20%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06595'
$ws.Range('E18').Value = '  -0.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '83.33'
$ws.Range('E19').Value = '  -2.70%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.000'
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.135'
$ws.Range('E21').Value = '  -5.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '15.37'
$ws.Range('E22').Value = '  -4.79%  '
$ws.Range('E23').Value = '  -7.65%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.366'
$ws.Range('E24').Value = '  +4.46%  '
$ws.Range('D25').Value = '21.791.07'
$ws.Range('E25').Value = '  -1.65%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.387'
$ws.Range('E26').Value = '  -6.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '145.38'
$ws.Range('E27').Value = '  -4.00%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.41'
$ws.Range('E28').Value = '  -3.66%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.843'
$ws.Range('E29').Value = '  -0.29%  '
$ws.Range('D30').Value = '1.715.47'
$ws.Range('E30').Value = '  -1.13%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '117.72'
$ws.Range('E31').Value = '  -3.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.913'
$ws.Range('E32').Value = '  -1.95%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.9640'
$ws.Range('E33').Value = '  -14.80%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08197'
$ws.Range('E34').Value = '  -0.37%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '8.966'
$ws.Range('E35').Value = '  -4.32%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06129'
$ws.Range('E36').Value = '  -1.50%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.125'
$ws.Range('E37').Value = '  -3.46%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02212'
$ws.Range('E38').Value = '  -4.20%  '
$ws.Range('B39').Value = 'WEMIXTOKEN'
$ws.Range('C39').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.454'
$ws.Range('E39').Value = '  -22.39%  '
$ws.Range('B40').Value = 'Algorand'
$ws.Range('C40').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2043'
$ws.Range('E40').Value = '  -4.73%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.183'
$ws.Range('E41').Value = '  -4.20%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9996'
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('E43').Value = '  -3.66%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5758'
$ws.Range('E44').Value = '  -3.88%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.13'
$ws.Range('E45').Value = '  -4.14%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.739'
$ws.Range('E46').Value = '  -0.55%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5527'
$ws.Range('E47').Value = '  -4.63%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '117.99'
$ws.Range('E48').Value = '  -2.79%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.875'
$ws.Range('E49').Value = '  -5.79%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.136'
$ws.Range('E50').Value = '  -3.17%  '
